$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '65.561.28'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.17%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.650.78'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.79%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.80%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.01'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.99%  '

# Row 7
$ws.Range('B7').Value = 'USDC'
$ws.Range('C7').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.640'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.33%  '

# Row 9
$ws.Range('E9').Value = '  -2.54%  '

# Row 10
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.68%  '

# Row 11
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.396'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.50%  '

# Row 12
$ws.Range('E12').Value = '  +1.05%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.94'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.88%  '

# Row 14
$ws.Range('E14').Value = '  -1.51%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.128.99'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.72%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.459.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.10%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.685.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.38%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.45'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.22%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.75'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.91%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '352.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.09%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.42'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.95%  '

# Row 22
$ws.Range('E22').Value = '  -0.04%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '69.42'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.41%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000112'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.68%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.75'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.73%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.62'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.04%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.61'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.06%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '556.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.28%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.05'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.16%  '

# Row 30
$ws.Range('E30').Value = '  -2.37%  '

# Row 31
$ws.Range('E31').Value = '  +0.31%  '

# Row 32
$ws.Range('E32').Value = '  -0.98%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.79'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.65%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.63'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.12%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.44'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.58%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.420'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.60%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '20.41'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.58%  '

# Row 38
$ws.Range('E38').Value = '  +0.08%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.43%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '153.48'
$ws.Range('D40').Style = 'Normal'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.01%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.59%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '160.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.60%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.41%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0613'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.72%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.11'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.65%  '

# Row 47
$ws.Range('E47').Value = '  +0.60%  '

# Row 48
$ws.Range('E48').Value = '  -0.20%  '

# Row 49
$ws.Range('E49').Value = '  +1.84%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.44%  '

# Row 51
$ws.Range('E51').Value = '  -8.04%  '
